$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 568 (shifts existing rows 568-601 down to 569-602)
$ws.Rows.Item(568).Insert()

# Populate the newly inserted row 568 with the new data record
$ws.Range("A568").Value = 9
$ws.Range("B568").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C568").Value = "Metropolitana"
$ws.Range("D568").Value = 44826
$ws.Range("E568").Value = 13
$ws.Range("F568").Value = 100112031
$ws.Range("G568").Value = "Poroto verde"
$ws.Range("H568").Value = "Magnum"
$ws.Range("I568").Value = "Primera"
$ws.Range("J568").Value = 65
$ws.Range("K568").Value = 29000
$ws.Range("L568").Value = 30000
$ws.Range("M568").Value = 29462
$ws.Range("N568").Value = "$/malla 25 kilos"
$ws.Range("O568").Value = "Perú"
$ws.Range("P568").Value = 1178
$ws.Range("Q568").Value = 25
$ws.Range("R568").Value = "Hortaliza"
